$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.362.28"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  -0.27%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'1.847.92"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  -0.10%  "
$ws.Range('E3').ClearFormats()
$ws.Range('D4').Value = "'0.9986"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "'  -0.15%  "
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'240.35"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  -0.15%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'0.6287"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  -0.14%  "
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'0.9992"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  -0.12%  "
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'0.07588"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'  -1.48%  "
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'0.2918"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  -0.42%  "
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'24.51"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  -0.74%  "
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'0.07743"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  +0.02%  "
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'1.848.25"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  -1.75%  "
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'5.010"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  -0.49%  "
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = "'0.6781"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  -0.11%  "
$ws.Range('E14').ClearFormats()
$ws.Range('E15').Value = "'  -3.25%  "
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'83.14"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  -0.61%  "
$ws.Range('E16').ClearFormats()
$ws.Range('B17').Value = "'Uniswap"
$ws.Range('B17').ClearFormats()
$ws.Range('C17').Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range('C17').ClearFormats()
$ws.Range('D17').Value = "'6.112"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'  -1.35%  "
$ws.Range('E17').ClearFormats()
$ws.Range('B18').Value = "'WrappedBTC"
$ws.Range('B18').ClearFormats()
$ws.Range('C18').Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('C18').ClearFormats()
$ws.Range('D18').Value = "'29.361.17"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  -0.38%  "
$ws.Range('E18').ClearFormats()
$ws.Range('B19').Value = "'BitcoinCash"
$ws.Range('B19').ClearFormats()
$ws.Range('C19').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('C19').ClearFormats()
$ws.Range('D19').Value = "'229.77"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  +0.66%  "
$ws.Range('E19').ClearFormats()
$ws.Range('B20').Value = "'Avalanche"
$ws.Range('B20').ClearFormats()
$ws.Range('C20').Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('C20').ClearFormats()
$ws.Range('D20').Value = "'12.35"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "'  -0.83%  "
$ws.Range('E20').ClearFormats()
$ws.Range('B21').Value = "'Dai"
$ws.Range('B21').ClearFormats()
$ws.Range('C21').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('C21').ClearFormats()
$ws.Range('D21').Value = "'0.9995"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  -0.04%  "
$ws.Range('E21').ClearFormats()
$ws.Range('B22').Value = "'Chainlink"
$ws.Range('B22').ClearFormats()
$ws.Range('C22').Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('C22').ClearFormats()
$ws.Range('D22').Value = "'7.428"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  -0.11%  "
$ws.Range('E22').ClearFormats()
$ws.Range('B23').Value = "'BinanceUSD"
$ws.Range('B23').ClearFormats()
$ws.Range('C23').Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range('C23').ClearFormats()
$ws.Range('D23').Value = "'0.9993"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  -0.12%  "
$ws.Range('E23').ClearFormats()
$ws.Range('B24').Value = "'Monero"
$ws.Range('B24').ClearFormats()
$ws.Range('C24').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C24').ClearFormats()
$ws.Range('D24').Value = "'159.00"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  +0.80%  "
$ws.Range('E24').ClearFormats()
$ws.Range('B25').Value = "'Stellar"
$ws.Range('B25').ClearFormats()
$ws.Range('C25').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('C25').ClearFormats()
$ws.Range('D25').Value = "'0.1395"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  +1.13%  "
$ws.Range('E25').ClearFormats()
$ws.Range('B26').Value = "'Cosmos"
$ws.Range('B26').ClearFormats()
$ws.Range('C26').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('C26').ClearFormats()
$ws.Range('D26').Value = "'8.442"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'  +0.30%  "
$ws.Range('E26').ClearFormats()
$ws.Range('B27').Value = "'EthereumClassic"
$ws.Range('B27').ClearFormats()
$ws.Range('C27').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C27').ClearFormats()
$ws.Range('D27').Value = "'17.64"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  -0.31%  "
$ws.Range('E27').ClearFormats()
$ws.Range('B28').Value = "'Toncoin"
$ws.Range('B28').ClearFormats()
$ws.Range('C28').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('C28').ClearFormats()
$ws.Range('D28').Value = "'1.430"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  +6.36%  "
$ws.Range('E28').ClearFormats()
$ws.Range('B29').Value = "'PancakeSwap"
$ws.Range('B29').ClearFormats()
$ws.Range('C29').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C29').ClearFormats()
$ws.Range('D29').Value = "'1.473"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  +0.46%  "
$ws.Range('E29').ClearFormats()
$ws.Range('B30').Value = "'Hedera"
$ws.Range('B30').ClearFormats()
$ws.Range('C30').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('C30').ClearFormats()
$ws.Range('D30').Value = "'0.05674"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  -0.08%  "
$ws.Range('E30').ClearFormats()
$ws.Range('B31').Value = "'Filecoin"
$ws.Range('B31').ClearFormats()
$ws.Range('C31').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C31').ClearFormats()
$ws.Range('D31').Value = "'4.116"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  -0.26%  "
$ws.Range('E31').ClearFormats()
$ws.Range('B32').Value = "'InternetComputer(DFINITY)"
$ws.Range('B32').ClearFormats()
$ws.Range('C32').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('C32').ClearFormats()
$ws.Range('D32').Value = "'4.037"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'  +0.02%  "
$ws.Range('E32').ClearFormats()
$ws.Range('B33').Value = "'ARBITRUM"
$ws.Range('B33').ClearFormats()
$ws.Range('C33').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('C33').ClearFormats()
$ws.Range('D33').Value = "'1.155"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  -0.71%  "
$ws.Range('E33').ClearFormats()
$ws.Range('B34').Value = "'LidoDAOToken"
$ws.Range('B34').ClearFormats()
$ws.Range('C34').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('C34').ClearFormats()
$ws.Range('D34').Value = "'1.822"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  -1.32%  "
$ws.Range('E34').ClearFormats()
$ws.Range('B35').Value = "'ImmutableX"
$ws.Range('B35').ClearFormats()
$ws.Range('C35').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('C35').ClearFormats()
$ws.Range('D35').Value = "'0.6972"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'  -0.90%  "
$ws.Range('E35').ClearFormats()
$ws.Range('B36').Value = "'HuobiToken"
$ws.Range('B36').ClearFormats()
$ws.Range('C36').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('C36').ClearFormats()
$ws.Range('D36').Value = "'2.579"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'  -0.24%  "
$ws.Range('E36').ClearFormats()
$ws.Range('B37').Value = "'VeChain"
$ws.Range('B37').ClearFormats()
$ws.Range('C37').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C37').ClearFormats()
$ws.Range('D37').Value = "'0.01824"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  +1.81%  "
$ws.Range('E37').ClearFormats()
$ws.Range('B38').Value = "'Maker"
$ws.Range('B38').ClearFormats()
$ws.Range('C38').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('C38').ClearFormats()
$ws.Range('D38').Value = "'1.236.87"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'  +1.34%  "
$ws.Range('E38').ClearFormats()
$ws.Range('B39').Value = "'MXToken"
$ws.Range('B39').ClearFormats()
$ws.Range('C39').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('C39').ClearFormats()
$ws.Range('D39').Value = "'2.715"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  -2.34%  "
$ws.Range('E39').ClearFormats()
$ws.Range('B40').Value = "'FraxShare"
$ws.Range('B40').ClearFormats()
$ws.Range('C40').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('C40').ClearFormats()
$ws.Range('D40').Value = "'6.424"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  -1.95%  "
$ws.Range('E40').ClearFormats()
$ws.Range('B41').Value = "'TrustWalletToken"
$ws.Range('B41').ClearFormats()
$ws.Range('C41').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('C41').ClearFormats()
$ws.Range('D41').Value = "'0.9017"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  -0.35%  "
$ws.Range('E41').ClearFormats()
$ws.Range('B42').Value = "'PaxDollar"
$ws.Range('B42').ClearFormats()
$ws.Range('C42').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('C42').ClearFormats()
$ws.Range('D42').Value = "'0.9991"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'  -0.17%  "
$ws.Range('E42').ClearFormats()
$ws.Range('B43').Value = "'RocketPoolETH"
$ws.Range('B43').ClearFormats()
$ws.Range('C43').Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range('C43').ClearFormats()
$ws.Range('D43').Value = "'2.007.90"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'  -2.20%  "
$ws.Range('E43').ClearFormats()
$ws.Range('B44').Value = "'Quant"
$ws.Range('B44').ClearFormats()
$ws.Range('C44').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('C44').ClearFormats()
$ws.Range('D44').Value = "'101.38"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  -0.49%  "
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = "'Aave"
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = "'65.70"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -0.87%  "
$ws.Range('E45').ClearFormats()
$ws.Range('B46').Value = "'Aptos"
$ws.Range('B46').ClearFormats()
$ws.Range('C46').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('C46').ClearFormats()
$ws.Range('D46').Value = "'7.135"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'  -0.02%  "
$ws.Range('E46').ClearFormats()
$ws.Range('B47').Value = "'BabyDogeCoin"
$ws.Range('B47').ClearFormats()
$ws.Range('C47').Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range('C47').ClearFormats()
$ws.Range('D47').Value = "'0.00000000118"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  -2.21%  "
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'0.1156"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  +0.91%  "
$ws.Range('E48').ClearFormats()
$ws.Range('B49').Value = "'EnergySwap"
$ws.Range('B49').ClearFormats()
$ws.Range('C49').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C49').ClearFormats()
$ws.Range('D49').Value = "'9.029"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  +0.20%  "
$ws.Range('E49').ClearFormats()
$ws.Range('B50').Value = "'TheSandbox"
$ws.Range('B50').ClearFormats()
$ws.Range('C50').Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('C50').ClearFormats()
$ws.Range('D50').Value = "'0.3978"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  -1.03%  "
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = "'1.676"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'  -0.35%  "
$ws.Range('E51').ClearFormats()
